$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the three runs "Storage (Week 8" / ", 12" / ")" into a
# single run reading "Storage (Week 8, 12)". A find/replace over the whole
# phrase (which spans the three original runs) collapses them into one run
# because the runs all share identical formatting.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Storage (Week 8, 12)", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "Storage (Week 8, 12)", 2)

# ---------------------------------------------------------------------------
# Change 2: split "Computations (plus, minus, modulus, d" / "ivide,
# multiplication,.)" into six runs, breaking exactly on the word boundaries
# that Word's proofing pass would flag ("modulus" and "d" + "ivide"):
#   "Computations (plus, minus, " | "modulus" | ", " | "d" | "ivide" |
#   ", multiplication,.)"
# Toggling a character-formatting property on a sub-range and then toggling
# it back forces the engine to split the run at those boundaries while
# leaving the effective formatting (and so the serialized <w:rPr>) unchanged.
# ---------------------------------------------------------------------------
$compPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Computations (plus, minus, modulus*") {
        $compPara = $cand
        break
    }
}

$compStart = $compPara.Range.Start
$paraTextEnd = $compPara.Range.End - 1   # exclude the paragraph mark

$rMod = $d.Range($compStart, $paraTextEnd)
$rMod.Find.Execute("modulus")
$modStart = $rMod.Start
$modEnd = $rMod.End

$rDiv = $d.Range($modEnd, $paraTextEnd)
$rDiv.Find.Execute("divide")
$divStart = $rDiv.Start
$divEnd = $rDiv.End

$segments = @(
    @($compStart, $modStart),
    @($modStart, $modEnd),
    @($modEnd, $divStart),
    @($divStart, $divStart + 1),
    @($divStart + 1, $divEnd),
    @($divEnd, $paraTextEnd)
)

foreach ($seg in $segments) {
    $sr = $d.Range($seg[0], $seg[1])
    $sr.Font.Bold = $true
    $sr.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# Change 3: remove the "Namespaces (brief)" paragraph that immediately
# follows the "References" paragraph (a different, earlier "Namespaces
# (brief)" paragraph elsewhere in the document is left untouched).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $refPara = $d.Paragraphs.Item($i)
    if ($refPara.Range.Text -eq "References`r") {
        $nextPara = $d.Paragraphs.Item($i + 1)
        if ($nextPara.Range.Text -eq "Namespaces (brief)`r") {
            $nextPara.Range.Delete()
        }
        break
    }
}
